$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 195.32
$ws.Range("I15").Value = 195.32
$ws.Range("K15").Value = 585.96
$ws.Range("M15").Value = -416.96
$ws.Range("H28").Value = 520.7692
$ws.Range("I28").Value = 303
$ws.Range("J28").Value = 869.2
$ws.Range("K28").Value = 303
$ws.Range("L28").Value = 869.2
$ws.Range("M28").Value = 182
$ws.Range("N28").Value = -1839.2
$ws.Range("H86").Value = 4064.5417
$ws.Range("J86").Value = 2697.75
$ws.Range("L86").Value = 2697.75
$ws.Range("N86").Value = -4943.75
$ws.Range("H88").Value = 25499.4
$ws.Range("I88").Value = 1500
$ws.Range("J88").Value = 31499.25
$ws.Range("K88").Value = 1500
$ws.Range("L88").Value = 31499.25
$ws.Range("M88").Value = -1094
$ws.Range("N88").Value = -32311.25
$ws.Range("H89").Value = 4064.5417
$ws.Range("J89").Value = 2697.75
$ws.Range("L89").Value = 13488.75
$ws.Range("N89").Value = -24720.75
$ws.Range("H91").Value = 25499.4
$ws.Range("I91").Value = 1500
$ws.Range("J91").Value = 31499.25
$ws.Range("K91").Value = 1500
$ws.Range("L91").Value = 31499.25
$ws.Range("M91").Value = -96
$ws.Range("N91").Value = -34307.25
$ws.Range("H111").Value = 1325
$ws.Range("I111").Value = 1325
$ws.Range("K111").Value = 3975
$ws.Range("M111").Value = -908
$ws.Range("H112").Value = 18434676
$ws.Range("I112").Value = 650
$ws.Range("J112").Value = 21979682
$ws.Range("K112").Value = 1950
$ws.Range("L112").Value = 65939046
$ws.Range("M112").Value = -842
$ws.Range("N112").Value = -65941262
$ws.Range("H113").Value = 20835600
$ws.Range("I113").Value = 2600
$ws.Range("J113").Value = 41668600
$ws.Range("K113").Value = 2600
$ws.Range("L113").Value = 41668600
$ws.Range("M113").Value = 654
$ws.Range("N113").Value = -41675108
$ws.Range("H137").Value = 1222.3334
$ws.Range("I137").Value = 1216.3928
$ws.Range("K137").Value = 3649.1784
$ws.Range("M137").Value = -1099.1784

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 154684.77
$ws.Range("I2").Value = 167325.17
$ws.Range("J2").Value = 3000
$ws.Range("K2").Value = 167325.17
$ws.Range("L2").Value = 3000
$ws.Range("M2").Value = -167212.17
$ws.Range("N2").Value = -3226
$ws.Range("H32").Value = 15935.949
$ws.Range("I32").Value = 11792.743
$ws.Range("J32").Value = 28199.84
$ws.Range("K32").Value = 11792.743
$ws.Range("L32").Value = 28199.84
$ws.Range("M32").Value = -11505.743
$ws.Range("N32").Value = -28773.84
$ws.Range("H74").Value = 10205461
$ws.Range("I74").Value = 1131.9697
$ws.Range("J74").Value = 31251892
$ws.Range("K74").Value = 1131.9697
$ws.Range("L74").Value = 31251892
$ws.Range("M74").Value = -257.9697000000001
$ws.Range("N74").Value = -31253640
$ws.Range("H77").Value = 10205461
$ws.Range("I77").Value = 1131.9697
$ws.Range("J77").Value = 31251892
$ws.Range("K77").Value = 5659.8485
$ws.Range("L77").Value = 156259460
$ws.Range("M77").Value = -1291.8485
$ws.Range("N77").Value = -156268196
$ws.Range("H109").Value = 41369.25
$ws.Range("J109").Value = 41369.25
$ws.Range("L109").Value = 41369.25
$ws.Range("N109").Value = -44143.25
$ws.Range("H116").Value = 154684.77
$ws.Range("I116").Value = 167325.17
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 167325.17
$ws.Range("L116").Value = 3000
$ws.Range("M116").Value = -165031.17
$ws.Range("N116").Value = -7588
$ws.Range("H133").Value = 33940
$ws.Range("J133").Value = 33940
$ws.Range("L133").Value = 33940
$ws.Range("N133").Value = -39000

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 154684.77
$ws.Range("I3").Value = 167325.17
$ws.Range("J3").Value = 3000
$ws.Range("K3").Value = 167325.17
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = -167211.17
$ws.Range("N3").Value = -3228
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 1640.6
$ws.Range("I12").Value = 1640.6
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 1640.6
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -1470.6
$ws.Range("N12").ClearContents()
$ws.Range("H31").Value = 6949372
$ws.Range("I31").Value = 1469.5758
$ws.Range("J31").Value = 12828366
$ws.Range("K31").Value = 1469.5758
$ws.Range("L31").Value = 12828366
$ws.Range("M31").Value = -1174.5758
$ws.Range("N31").Value = -12828956
$ws.Range("H34").Value = 6949372
$ws.Range("I34").Value = 1469.5758
$ws.Range("J34").Value = 12828366
$ws.Range("K34").Value = 1469.5758
$ws.Range("L34").Value = 12828366
$ws.Range("M34").Value = -1267.5758
$ws.Range("N34").Value = -12828770
$ws.Range("H43").Value = 41328.5
$ws.Range("J43").Value = 41328.5
$ws.Range("L43").Value = 41328.5
$ws.Range("N43").Value = -41696.5
$ws.Range("H99").Value = 25003600
$ws.Range("I99").Value = 3000
$ws.Range("J99").Value = 31253750
$ws.Range("K99").Value = 3000
$ws.Range("L99").Value = 31253750
$ws.Range("M99").Value = -1502
$ws.Range("N99").Value = -31256746
$ws.Range("H101").Value = 41328.5
$ws.Range("J101").Value = 41328.5
$ws.Range("L101").Value = 41328.5
$ws.Range("N101").Value = -47818.5
$ws.Range("H126").Value = 25003600
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 31253750
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 93761250
$ws.Range("M126").Value = -6530
$ws.Range("N126").Value = -93766190

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 87.666664
$ws.Range("I11").Value = 87.666664
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 262.999992
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -122.999992
$ws.Range("N11").ClearContents()
$ws.Range("H81").Value = 5267.143
$ws.Range("I81").Value = 1548
$ws.Range("J81").Value = 7333.3335
$ws.Range("K81").Value = 4644
$ws.Range("L81").Value = 22000.0005
$ws.Range("M81").Value = -3521
$ws.Range("N81").Value = -24246.0005
$ws.Range("H82").Value = 1003
$ws.Range("I82").Value = 1003
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 3009
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -2603
$ws.Range("N82").ClearContents()
$ws.Range("H84").Value = 5267.143
$ws.Range("I84").Value = 1548
$ws.Range("J84").Value = 7333.3335
$ws.Range("K84").Value = 13932
$ws.Range("L84").Value = 66000.0015
$ws.Range("M84").Value = -8316
$ws.Range("N84").Value = -77232.0015
$ws.Range("H85").Value = 1003
$ws.Range("I85").Value = 1003
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 3009
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -1605
$ws.Range("N85").ClearContents()
$ws.Range("H114").Value = 5188.393
$ws.Range("I114").Value = 298
$ws.Range("J114").Value = 7144.55
$ws.Range("K114").Value = 894
$ws.Range("L114").Value = 21433.65
$ws.Range("M114").Value = 2360
$ws.Range("N114").Value = -27941.65
$ws.Range("H129").Value = 1246.5217
$ws.Range("I129").Value = 1035
$ws.Range("J129").Value = 1477.2727
$ws.Range("K129").Value = 3105
$ws.Range("L129").Value = 4431.8181
$ws.Range("M129").Value = 1895
$ws.Range("N129").Value = -14431.8181
$ws.Range("H132").Value = 5558145.5
$ws.Range("J132").Value = 10104465
$ws.Range("L132").Value = 90940185
$ws.Range("N132").Value = -90945245

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3175.818
$ws.Range("I102").Value = 2557.6667
$ws.Range("J102").Value = 3917.6
$ws.Range("K102").Value = 2557.6667
$ws.Range("L102").Value = 3917.6
$ws.Range("M102").Value = -935.6667000000002
$ws.Range("N102").Value = -7161.6
$ws.Range("H132").Value = 7938755
$ws.Range("I132").Value = 11906548
$ws.Range("J132").Value = 3170.5715
$ws.Range("K132").Value = 35719644
$ws.Range("L132").Value = 9511.7145
$ws.Range("M132").Value = -35717114
$ws.Range("N132").Value = -14571.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 2757.1428
$ws.Range("J38").Value = 2757.1428
$ws.Range("L38").Value = 2757.1428
$ws.Range("N38").Value = -3577.1428
$ws.Range("H103").Value = 36801
$ws.Range("J103").Value = 36801
$ws.Range("L103").Value = 36801
$ws.Range("N103").Value = -39145
$ws.Range("H133").Value = 80108.664
$ws.Range("J133").Value = 80108.664
$ws.Range("L133").Value = 80108.664
$ws.Range("N133").Value = -85168.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H28").Value = 7216.3335
$ws.Range("J28").Value = 7216.3335
$ws.Range("L28").Value = 7216.3335
$ws.Range("N28").Value = -7912.3335
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
$ws.Range("H33").Value = 5000
$ws.Range("I33").Value = 5000
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 5000
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -4750
$ws.Range("N33").ClearContents()
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H36").Value = 5000
$ws.Range("I36").Value = 5000
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 5000
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -4750
$ws.Range("N36").ClearContents()
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H126").Value = 794.25
$ws.Range("I126").Value = 729.9545000000001
$ws.Range("K126").Value = 2189.8635
$ws.Range("M126").Value = 280.1364999999996
$ws.Range("H136").Value = 9808646
$ws.Range("I136").Value = 4679.2856
$ws.Range("K136").Value = 14037.8568
$ws.Range("M136").Value = -11487.8568
